# Weekly update: a new price record for "Macroferia Regional de Talca - Alcachofa"
# was reported. It is inserted as a new row right before the previous most
# recent entry (row 39), pushing that record and all the ones below it down
# by one row. The new record duplicates the previous week's values, only the
# date (column D) changes to the new reporting date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 39, shifting existing rows 39-49 down to 40-50.
$ws.Rows.Item(39).Insert()

# Populate the new row 39 with the new weekly record.
$ws.Range("A39").Value = 5
$ws.Range("B39").Value = "Macroferia Regional de Talca"
$ws.Range("C39").Value = "Maule"
$ws.Range("D39").Value = 44455
$ws.Range("E39").Value = 7
$ws.Range("F39").Value = 100112013
$ws.Range("G39").Value = "Alcachofa"
$ws.Range("H39").Value = "Madrigal"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 300
$ws.Range("K39").Value = 12000
$ws.Range("L39").Value = 12000
$ws.Range("M39").Value = 12000
$ws.Range("N39").Value = "$/caja 40 unidades"
$ws.Range("O39").Value = "Provincia del Elquí"
$ws.Range("P39").Value = 300
$ws.Range("Q39").Value = 40
$ws.Range("R39").Value = "Hortaliza"
